$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add P1=14, Q1=15 and give them the same style as the
# neighboring header cell (O1), by copying its formatting rather than rebuilding
# the style from scratch (keeps styles.xml untouched, reusing the existing xf). ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Data rows 2-25 ---
# Existing columns I, K, M, O swap their 1/2 values; new columns P and Q are
# added (unstyled, like the other data columns) with value 2.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value  = 2   # I: was 1 -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K: was 2 -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M: was 1 -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O: was 2 -> 1
    $ws.Cells.Item($r, 16).Value = 2   # P: new
    $ws.Cells.Item($r, 17).Value = 2   # Q: new
}
